$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-tracking entry for 29.11.2019, "IO Control Unit" / "Comments"
$ws.Range("A11").Value = "29.11.2019"
$ws.Range("B11").Value = 0.64583333333333337
$ws.Range("C11").Value = 0.66666666666666663
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("E11").Value = "IO Control Unit"
$ws.Range("F11").Value = "Comments"

# Match number formats used by the rows above (B/C = built-in h:mm, D = custom [$]hh:mm;@)
$ws.Range("B11:C11").NumberFormat = "h:mm"
$ws.Range("D11").NumberFormat = "[$]hh:mm;@"

# Move the active selection the way Excel leaves it after typing the new row
$ws.Range("A12").Select()
